$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in A4 (next row after existing data)
$ws.Range("A4").Value = 4

# Move the active selection to A5, mirroring the original "next empty cell" selection
$ws.Range("A5").Select()
